# Froze 2 decoder layer
# Update ASR Results sheet: refresh column B (predicted) labels and column C (count) values
# for several rows to reflect the latest decoder run.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column B text corrections (predicted token changed) ---
$ws.Cells.Item(7, 2).Value = "<and>"
$ws.Cells.Item(9, 2).Value = "<zulu>"
$ws.Cells.Item(18, 2).Value = "<of>"
$ws.Cells.Item(30, 2).Value = "<to>"
$ws.Cells.Item(33, 2).Value = "<by>"
$ws.Cells.Item(46, 2).Value = "<there>"
$ws.Cells.Item(48, 2).Value = "<up>"
$ws.Cells.Item(50, 2).Value = "<xir>"

# --- Column C numeric updates (row => new value) ---
$updates = @{
    2  = 9
    3  = 5
    4  = 3
    5  = 8
    8  = 8
    9  = 11
    10 = 3
    11 = 12
    13 = 8
    14 = 2
    16 = 4
    17 = 8
    18 = 3
    19 = 7
    20 = 8
    22 = 4
    23 = 5
    24 = 6
    25 = 9
    26 = 4
    27 = 9
    28 = 14
    29 = 6
    30 = 12
    31 = 7
    32 = 6
    34 = 9
    35 = 5
    36 = 6
    37 = 13
    38 = 7
    39 = 5
    40 = 5
    41 = 6
    42 = 4
    43 = 10
    44 = 5
    45 = 13
    46 = 9
    47 = 7
    48 = 11
    49 = 5
    50 = 10
    51 = 5
    52 = 4
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 3).Value = $updates[$row]
}
